# Apply 2022-12-19 crime data updates (column I = 2022 totals) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7076
$ws.Range("I3").Value = 7294
$ws.Range("I4").Value = 1679
$ws.Range("I5").Value = 690
$ws.Range("I6").Value = 8663
$ws.Range("I7").Value = 25402

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 233
$ws.Range("I3").Value = 256
$ws.Range("I6").Value = 239
$ws.Range("I7").Value = 787

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 188
$ws.Range("I7").Value = 591

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 205
$ws.Range("I7").Value = 796
$ws.Range("I8").Value = 1507
$ws.Range("I11").Value = 389
$ws.Range("I19").Value = 708
$ws.Range("I27").Value = 221
$ws.Range("I29").Value = 1510
$ws.Range("I31").Value = 255
$ws.Range("I32").Value = 30
$ws.Range("I33").Value = 1119
$ws.Range("I36").Value = 350
$ws.Range("I37").Value = 787
$ws.Range("I40").Value = 46
$ws.Range("I42").Value = 956
$ws.Range("I48").Value = 323
$ws.Range("I51").Value = 294
$ws.Range("I63").Value = 76
$ws.Range("I65").Value = 591
$ws.Range("I76").Value = 364
$ws.Range("I77").Value = 155
$ws.Range("I78").Value = 337
$ws.Range("I79").Value = 731
$ws.Range("I83").Value = 550
$ws.Range("I85").Value = 1131
$ws.Range("I87").Value = 68
$ws.Range("I91").Value = 269
$ws.Range("I95").Value = 391
$ws.Range("I101").Value = 25402

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 201
$ws.Range("I4").Value = 23
$ws.Range("I7").Value = 550

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 254
$ws.Range("I7").Value = 1119

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 446
$ws.Range("I3").Value = 514
$ws.Range("I6").Value = 420
$ws.Range("I7").Value = 1510

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 206
$ws.Range("I7").Value = 708

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I4").Value = 43
$ws.Range("I6").Value = 164
$ws.Range("I7").Value = 323

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 168
$ws.Range("I7").Value = 364

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 324
$ws.Range("I3").Value = 425
$ws.Range("I7").Value = 1131

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 386
$ws.Range("I7").Value = 956

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 337

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 269

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 216
$ws.Range("I7").Value = 731

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 101
$ws.Range("I7").Value = 350

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 110
$ws.Range("I7").Value = 389

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 78
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 435
$ws.Range("I4").Value = 93
$ws.Range("I6").Value = 488
$ws.Range("I7").Value = 1507

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 221

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 120
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 217
$ws.Range("I7").Value = 796

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 68
